# Update (Removed Auto Arima)
# Recompute the Prophet / Amazon forecast columns on the "Forecast Comparison"
# sheet and propagate the updated totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# New values for columns C (Prophet Forecast), D (Amazon Mean Forecast),
# E (Amazon P70 Forecast), F (Amazon P80 Forecast), G (Amazon P90 Forecast)
# for rows 2-17 (weeks W01-W16).
$data = @(
    @(166, 189, 222, 252, 298),
    @(164, 185, 221, 258, 314),
    @(156, 146, 174, 204, 249),
    @(143, 121, 144, 168, 205),
    @(132, 121, 145, 171, 211),
    @(126, 120, 144, 169, 208),
    @(126, 118, 143, 170, 212),
    @(128, 120, 145, 172, 216),
    @(133, 117, 141, 165, 202),
    @(138, 117, 141, 167, 209),
    @(147, 118, 143, 171, 215),
    @(159, 120, 145, 174, 221),
    @(173, 122, 148, 178, 226),
    @(187, 113, 138, 167, 215),
    @(194, 116, 141, 171, 220),
    @(189, 111, 135, 163, 209)
)

$row = 2
foreach ($vals in $data) {
    $wsForecast.Cells.Item($row, 3).Value = $vals[0]  # C - Prophet Forecast
    $wsForecast.Cells.Item($row, 4).Value = $vals[1]  # D - Amazon Mean Forecast
    $wsForecast.Cells.Item($row, 5).Value = $vals[2]  # E - Amazon P70 Forecast
    $wsForecast.Cells.Item($row, 6).Value = $vals[3]  # F - Amazon P80 Forecast
    $wsForecast.Cells.Item($row, 7).Value = $vals[4]  # G - Amazon P90 Forecast
    $row++
}

# Update the Summary sheet totals derived from the Prophet Forecast column.
# These cells hold numeric-looking text (t="inlineStr"), so force the
# number format to Text before assigning, then restore the style so the
# cell keeps using the default "Normal" style (just like the original).
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $wsSummary.Range("B9")  "2461"   # Total Forecast (16 Weeks)
Set-TextValue $wsSummary.Range("B10") "1141"   # Total Forecast (8 Weeks)
Set-TextValue $wsSummary.Range("B11") "629"    # Total Forecast (4 Weeks)
Set-TextValue $wsSummary.Range("B12") "194"    # Max Forecast
Set-TextValue $wsSummary.Range("B14") "126"    # Min Forecast
